# Apply cryptocurrency price/volume updates scraped for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.317.37"
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("D3").Value = "2.529.43"
$ws.Range("E3").Value = "  -5.07%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'578.95"
$ws.Range("E5").Value = "  -2.74%  "
$ws.Range("D6").Value = "'169.28"
$ws.Range("E6").Value = "  -3.78%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("D9").Value = "2.529.01"
$ws.Range("E9").Value = "  -5.06%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'0.341"
$ws.Range("E12").Value = "  -4.43%  "
$ws.Range("D13").Value = "'4.87"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "2.991.13"
$ws.Range("E14").Value = "  -5.10%  "
$ws.Range("D15").Value = "70.180.67"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "'24.96"
$ws.Range("E17").Value = "  -4.56%  "
$ws.Range("D18").Value = "2.520.75"
$ws.Range("E18").Value = "  -5.19%  "
$ws.Range("D19").Value = "'11.41"
$ws.Range("E19").Value = "  -7.70%  "
$ws.Range("D20").Value = "'7.67"
$ws.Range("E20").Value = "  -5.62%  "
$ws.Range("D21").Value = "'354.74"
$ws.Range("E21").Value = "  -3.71%  "
$ws.Range("D22").Value = "'3.94"
$ws.Range("E22").Value = "  -6.25%  "
$ws.Range("D23").Value = "'1.99"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("D25").Value = "'69.24"
$ws.Range("E25").Value = "  -3.64%  "
$ws.Range("E26").Value = "  -5.62%  "
$ws.Range("D27").Value = "'9.16"
$ws.Range("E27").Value = "  -5.91%  "
$ws.Range("D28").Value = "2.658.74"
$ws.Range("E28").Value = "  -5.10%  "
$ws.Range("D29").Value = "'1.03"
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("D30").Value = "0.0₃0914"
$ws.Range("E30").Value = "  -5.07%  "
$ws.Range("D31").Value = "'7.85"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("D32").Value = "'484.76"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").Value = "'1.29"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "'1.77"
$ws.Range("E34").Value = "  -2.87%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.116"
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'155.28"
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("D38").Value = "'18.65"
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("D39").Value = "'18.89"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'4.78"
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").Value = "'1.63"
$ws.Range("E43").Value = "  -6.59%  "
$ws.Range("D44").Value = "'1.19"
$ws.Range("E44").Value = "  -13.15%  "
$ws.Range("E45").Value = "  -8.08%  "
$ws.Range("D46").Value = "'38.57"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").Value = "'143.75"
$ws.Range("E47").Value = "  -8.20%  "
$ws.Range("D48").Value = "'3.55"
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("D49").Value = "'0.529"
$ws.Range("E49").Value = "  -5.08%  "
$ws.Range("D50").Value = "'1.62"
$ws.Range("E50").Value = "  -5.62%  "
$ws.Range("D51").Value = "'0.599"
$ws.Range("E51").Value = "  -0.81%  "
